# Automatische test-sync: 2025-08-05 18:26:50
#
# Adds a new "Testmail #10" row to the Logs sheet, extends the
# conditional-formatting ranges to cover the new row, and fixes up the
# Dashboard category counts/order to reflect the new log entry.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 31 with the new test mail entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A31").Value = "Is er al nieuws?"
$logs.Range("B31").Value = "mailmind.test@zohomail.eu"
$logs.Range("C31").Value = "Testmail #10: Is er al nieuws?"
$logs.Range("D31").Value = "Opvolging / Status"
$logs.Range("E31").Value = "Dank voor je bericht. We hebben je eerdere e-mail ontvangen en doorgestuurd naar klantenservice@bedrijf.nl."
$logs.Range("F31").Value = "2025-08-05 18:26:02"
$logs.Range("G31").Value = "Ja"
$logs.Range("H31").Value = "Ja"
$logs.Range("I31").Value = "Nee"
$logs.Range("J31").Value = "Nee"

# ---------------------------------------------------------------------
# 2. Extend the conditional formatting ranges from row 30 to row 31
# ---------------------------------------------------------------------
$ranges = @("D2:D30", "G2:G30", "H2:H30", "I2:I30", "J2:J30")
foreach ($sq in $ranges) {
    $col = $sq.Substring(0, 1)
    $newRange = $logs.Range($col + "2:" + $col + "31")
    $fcs = $logs.Range($sq).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 3. Dashboard sheet: re-order category rows 5-7 and bump the count for
#    "Opvolging / Status" now that a second entry exists.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A5").Value = "Opvolging / Status"
$dash.Range("B5").Value = 2

$dash.Range("A6").Value = "Retour / Terugbetaling"
$dash.Range("B6").Value = 1

$dash.Range("A7").Value = "Klacht / Probleem"
$dash.Range("B7").Value = 1
